# Add two more logged time entries (rows 56 and 57) to the Time Log sheet,
# filling in Stop Time / Interruption / Activity for row 56 and the full
# row for row 57.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Row 56: Interruption first, then Stop Time, then Activity = Testing
$ws1.Range("D56").Value = 30
$ws1.Range("C56").Value = 0.83472222222222225
$ws1.Range("F56").Value = "Testing"

# Row 57: Interruption first, then Date/Start/Stop, then Activity = Coding
$ws1.Range("D57").Value = 0
$ws1.Range("A57").Value = 41902
$ws1.Range("B57").Value = 0.97361111111111109
$ws1.Range("C57").Value = 0.98749999999999993
$ws1.Range("F57").Value = "Coding"

# Update the active selection to reflect where the user ended up
$ws1.Range("B58").Select()

$wb.Save()
